$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 9 (pushes old row9..row30 down to row10..row31)
$ws.Rows("9:9").Insert()

# Insert a second new row at row 13 (after the first insert, this is the gap
# right after the rec_fr row, pushing everything below down one more row)
$ws.Rows("13:13").Insert()

# --- Fill in the new "A_heliostat" parameter row (row 9) ---
$ws.Range("B9").Value = "A_heliostat"
$ws.Range("C9").Value = 169
$ws.Range("D9").Value = "m2"
$ws.Range("E9").Value = "Heliostat reflective area"
$ws.Range("F9").Value = 0
$ws.Range("K9").Value = "assumed"

# --- Fill in the new "ab_rec" parameter row (row 13) ---
$ws.Range("B13").Value = "ab_rec"
$ws.Range("C13").Value = 0.96
$ws.Range("E13").Value = "Receiver coating absorptance"
$ws.Range("F13").Value = 0
$ws.Range("K13").Value = "assumed"

# --- Add the missing unit "h" to the t_storage row (now row 17) ---
$ws.Range("D17").Value = "h"

# Match the final selected cell left by the author's edit
[void]$ws.Range("J26").Select()
